$d = $word.ActiveDocument

# Remove the existing _GoBack bookmark; it will be re-created further down
# (on the newly inserted "Vi startede på 3 og sidste sprint..." paragraph),
# matching the diff where the bookmark moves to the end of the new content.
try {
    $d.Bookmarks.Item("_GoBack").Delete()
} catch {
}

# Locate the two paragraphs that need to be replaced:
#   "08-12-2017: Tiende dag"  and  "11-12-2017: Ellevte dag"
$count = $d.Paragraphs.Count
$startPara = $null
$endPara = $null
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "08-12-2017: Tiende dag") {
        $startPara = $p
    }
    if ($t -eq "11-12-2017: Ellevte dag") {
        $endPara = $p
    }
}

$rng = $d.Range($startPara.Range.Start, $endPara.Range.End)

$xml = '<w:p><w:r><w:t>08-12-2017</w:t></w:r><w:r><w:t>: Tiende dag</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Vi holdte i dag et møde med vores ’Product Owner’ Klaus, hvor vi fik vist vores produkt, hvordan det ser ud, hvad vi tænker og hvad han tænkte. Der blev skrevet noter ned om hans ønsker og efter mødet blev der arbejdet for at få gjort dette sprints </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>userstories</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> puttet over i ’DONE’. </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Vi sluttede med at alle </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>userstories</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> var opnået og endnu et vellykket sprint var overstået.</w:t></w:r></w:p><w:p><w:r><w:t>11-12-2017</w:t></w:r><w:r><w:t>: Ellevte dag</w:t></w:r></w:p><w:p><w:r><w:t>Vi startede på 3 og sidste sprint i dette projekt.</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/></w:p>'

$rng.InsertXML($xml)
